$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46052, 12222.5338800208, 11191.6786057924, 18795.86, 8242.84081945819, 26.6108093854405),
    @(3, 46053, 4965.32333927576, 7929.84174352245, 12075.86, 8111.91602651677, 165.245740418301),
    @(4, 46054, 5378.48890906186, 7961.64318850032, 9791.86, 8115.80238987394, 261.899399098928),
    @(5, 46055, 11849.0163547731, 11405.9690815929, 9791.86, 7927.95551457647, 397.586024840389),
    @(6, 46056, 11506.7705285637, 10905.8387779636, 9791.86, 7645.52522786993, 364.979333576396),
    @(7, 46057, 11506.7705285637, 10836.2427710317, 9791.86, 7645.52522786993, 362.079499954233),
    @(8, 46058, 11506.7705285637, 10910.349204469, 9791.86, 7645.52522786993, 365.16726801412),
    @(9, 46059, 11506.7705285637, 10130.3106172017, 9791.86, 7645.52522786993, 332.665660211318),
    @(10, 46060, 4877.58973915076, 7164.56084948764, 9791.86, 7276.68256210846, 193.724308816504),
    @(11, 46061, 4781.13067158898, 7168.02768543315, 9791.86, 7269.31477258062, 193.561769083907),
    @(12, 46062, 11119.1015365119, 10542.5018107652, 9791.86, 7175.53702002036, 330.257451282732),
    @(13, 46063, 11119.1015365119, 10745.6029501874, 9791.86, 7175.53702002036, 338.719998758658),
    @(14, 46064, 11119.1015365119, 10702.5265651126, 9791.86, 7175.53702002036, 336.92514938054),
    @(15, 46065, 11119.1015365119, 10804.440545741, 9791.86, 7181.99845348329, 341.440791634346)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
